$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E8 with the expanded note text (code sample appended) - this causes
# the old shared string at its old position to be removed (no longer referenced)
# and the new string appended at the end of the shared strings table, matching
# the observed index shuffle in the diff (52 removed, 53->52, 54->53, new->54).
$ws.Cells.Item(8, 5).Value = @'
This task can be completed by only giving the front end the full DNA sequence and an hash of exons positions/length; both retrievable with queries; all the front end needs is to know which sbstring within the main string to highligh; same we decided to do for the restriction sites task. Code for the front layer:                                                                 foreach my $key (keys %exons)   {
   substr($sequence, $key, $exons{$key}) = "<div>" . substr($sequence,    $key, $exons{$key} ) . "</div>";
}
print "<p> $sequence </p>";
'@

# Row 8 grows taller to fit the longer note.
$ws.Rows.Item(8).RowHeight = 165

# Move the active selection from F7 to F8 (sheetView selection change).
$ws.Cells.Item(8, 6).Select()

# Scroll the view down/right so C5 becomes the top-left visible cell
# (sheetView topLeftCell changes from C1 to C5).
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 3
